# Auto-generated edit script: updates specific profit-calculation cells
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3279.524
$ws.Range("I17").Value = 8563
$ws.Range("J17").Value = 2723.3684
$ws.Range("K17").Value = 25689
$ws.Range("L17").Value = 8170.1052
$ws.Range("M17").Value = -25521
$ws.Range("N17").Value = -8506.1052
# Row 40
$ws.Range("H40").Value = 2290.25
$ws.Range("I40").Value = 1797.6666
$ws.Range("J40").Value = 2585.8
$ws.Range("K40").Value = 1797.6666
$ws.Range("L40").Value = 2585.8
$ws.Range("M40").Value = -1622.6666
$ws.Range("N40").Value = -2935.8
# Row 112
$ws.Range("H112").Value = 1513.6364
$ws.Range("J112").Value = 1513.6364
$ws.Range("L112").Value = 4540.9092
$ws.Range("N112").Value = -6756.9092
# Row 129
$ws.Range("H129").Value = 866.2222
$ws.Range("J129").Value = 899.7857
$ws.Range("L129").Value = 2699.3571
$ws.Range("N129").Value = -12699.3571
# Row 132
$ws.Range("H132").Value = 1433.8
$ws.Range("I132").Value = 1443.8889
$ws.Range("K132").Value = 4331.6667
$ws.Range("M132").Value = -1801.6667
# Row 138
$ws.Range("H138").Value = 1911.806
$ws.Range("J138").Value = 2102.9697
$ws.Range("L138").Value = 6308.909100000001
$ws.Range("N138").Value = -16588.9091
# Row 140
$ws.Range("H140").Value = 63678.953
$ws.Range("J140").Value = 63678.953
$ws.Range("L140").Value = 63678.953
$ws.Range("N140").Value = -74038.95300000001
# Row 141
$ws.Range("H141").Value = 6488.2
$ws.Range("I141").Value = 2246.25
$ws.Range("K141").Value = 6738.75
$ws.Range("M141").Value = -1558.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4680.8423
$ws.Range("I32").Value = 3271.8367
$ws.Range("J32").Value = 13311
$ws.Range("K32").Value = 3271.8367
$ws.Range("L32").Value = 13311
$ws.Range("M32").Value = -2984.8367
$ws.Range("N32").Value = -13885
# Row 74
$ws.Range("H74").Value = 1753.5625
$ws.Range("I74").Value = 831.75
$ws.Range("J74").Value = 2675.375
$ws.Range("K74").Value = 831.75
$ws.Range("L74").Value = 2675.375
$ws.Range("M74").Value = 42.25
$ws.Range("N74").Value = -4423.375
# Row 77
$ws.Range("H77").Value = 1753.5625
$ws.Range("I77").Value = 831.75
$ws.Range("J77").Value = 2675.375
$ws.Range("K77").Value = 4158.75
$ws.Range("L77").Value = 13376.875
$ws.Range("M77").Value = 209.25
$ws.Range("N77").Value = -22112.875
# Row 132
$ws.Range("H132").Value = 1121.2084
$ws.Range("I132").Value = 1121.2084
$ws.Range("K132").Value = 3363.6252
$ws.Range("M132").Value = -833.6251999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 473.26086
$ws.Range("I107").Value = 437.55554
$ws.Range("J107").Value = 601.8
$ws.Range("K107").Value = 437.55554
$ws.Range("L107").Value = 601.8
$ws.Range("M107").Value = 1482.44446
$ws.Range("N107").Value = -4441.8
# Row 134
$ws.Range("H134").Value = 7453.0557
$ws.Range("I134").Value = 9480.615
$ws.Range("J134").Value = 2181.4
$ws.Range("K134").Value = 28441.845
$ws.Range("L134").Value = 6544.200000000001
$ws.Range("M134").Value = -25906.845
$ws.Range("N134").Value = -11614.2
# Row 137
$ws.Range("H137").Value = 61659.8
$ws.Range("J137").Value = 61659.8
$ws.Range("L137").Value = 61659.8
$ws.Range("N137").Value = -71859.8

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4832811
$ws.Range("I58").Value = 10870315
$ws.Range("J58").Value = 2808.2
$ws.Range("K58").Value = 10870315
$ws.Range("L58").Value = 2808.2
$ws.Range("N58").Value = -3214.2
$ws.Range("M58").Value = -10870112
# Row 70
$ws.Range("H70").Value = 46600
$ws.Range("J70").Value = 46600
$ws.Range("L70").Value = 46600
$ws.Range("N70").Value = -47230
# Row 73
$ws.Range("H73").Value = 46600
$ws.Range("J73").Value = 46600
$ws.Range("L73").Value = 46600
$ws.Range("N73").Value = -48784
# Row 99
$ws.Range("H99").Value = 3487.4285
$ws.Range("I99").Value = 2228
$ws.Range("K99").Value = 2228
$ws.Range("M99").Value = -730
# Row 126
$ws.Range("H126").Value = 3487.4285
$ws.Range("I126").Value = 2228
$ws.Range("K126").Value = 6684
$ws.Range("M126").Value = -4214
# Row 132
$ws.Range("H132").Value = 1728.4783
$ws.Range("I132").Value = 1429.625
$ws.Range("J132").Value = 2411.5715
$ws.Range("K132").Value = 4288.875
$ws.Range("L132").Value = 7234.7145
$ws.Range("M132").Value = -1758.875
$ws.Range("N132").Value = -12294.7145
# Row 134
$ws.Range("H134").Value = 2394.25
$ws.Range("I134").Value = 1938.1111
$ws.Range("J134").Value = 6499.5
$ws.Range("K134").Value = 5814.3333
$ws.Range("L134").Value = 19498.5
$ws.Range("M134").Value = -3279.3333
$ws.Range("N134").Value = -24568.5
# Row 136
$ws.Range("H136").Value = 4832811
$ws.Range("I136").Value = 10870315
$ws.Range("J136").Value = 2808.2
$ws.Range("K136").Value = 32610945
$ws.Range("L136").Value = 8424.599999999999
$ws.Range("N136").Value = -13524.6
$ws.Range("M136").Value = -32608395

$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 131
$ws.Range("H131").Value = 12281.549
$ws.Range("J131").Value = 12665.434
$ws.Range("L131").Value = 37996.302
$ws.Range("N131").Value = -48076.302
# Row 139
$ws.Range("H139").Value = 2061.2942
$ws.Range("I139").Value = 1979.3572
$ws.Range("J139").Value = 2443.6667
$ws.Range("K139").Value = 5938.071599999999
$ws.Range("L139").Value = 7331.000100000001
$ws.Range("M139").Value = -798.0715999999993
$ws.Range("N139").Value = -17611.0001
# Row 140
$ws.Range("H140").Value = 3463.625
$ws.Range("I140").Value = 1665.125
$ws.Range("K140").Value = 4995.375
$ws.Range("M140").Value = 184.625

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2175
$ws.Range("I97").Value = 2233.3333
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 2233.3333
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1737.3333
$ws.Range("N97").Value = -2992
# Row 122
$ws.Range("H122").Value = 1786.8572
$ws.Range("I122").Value = 1625
$ws.Range("J122").Value = 2002.6666
$ws.Range("K122").Value = 4875
$ws.Range("L122").Value = 6007.9998
$ws.Range("M122").Value = -2425
$ws.Range("N122").Value = -10907.9998
# Row 132
$ws.Range("H132").Value = 1751432.2
$ws.Range("I132").Value = 2406638.8
$ws.Range("J132").Value = 4215.1665
$ws.Range("K132").Value = 7219916.399999999
$ws.Range("L132").Value = 12645.4995
$ws.Range("M132").Value = -7217386.399999999
$ws.Range("N132").Value = -17705.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2184.0454
$ws.Range("I22").Value = 2118.0625
$ws.Range("K22").Value = 2118.0625
$ws.Range("M22").Value = -1823.0625
# Row 27
$ws.Range("H27").Value = 2184.0454
$ws.Range("I27").Value = 2118.0625
$ws.Range("K27").Value = 2118.0625
$ws.Range("M27").Value = -2011.0625
# Row 55
$ws.Range("H55").Value = 680.2727
$ws.Range("I55").Value = 667.1667
$ws.Range("K55").Value = 667.1667
$ws.Range("M55").Value = -494.1667
# Row 93
$ws.Range("H93").Value = 15152465
$ws.Range("I93").Value = 810.82355
$ws.Range("K93").Value = 810.82355
$ws.Range("M93").Value = 437.17645
# Row 100
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959
# Row 132
$ws.Range("H132").Value = 2062.0334
$ws.Range("I132").Value = 1224
$ws.Range("K132").Value = 3672
$ws.Range("M132").Value = -1142
# Row 136
$ws.Range("H136").Value = 4155.931
$ws.Range("I136").Value = 3364.125
$ws.Range("K136").Value = 10092.375
$ws.Range("M136").Value = -7542.375

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 108
$ws.Range("H108").Value = 67999.5
$ws.Range("J108").Value = 67999.5
$ws.Range("L108").Value = 67999.5
$ws.Range("N108").Value = -75679.5
# Row 122
$ws.Range("H122").Value = 56989.645
$ws.Range("J122").Value = 1779
$ws.Range("L122").Value = 5337
$ws.Range("N122").Value = -10237
# Row 132
$ws.Range("H132").Value = 1389.8064
$ws.Range("I132").Value = 1039.8572
$ws.Range("K132").Value = 3119.5716
$ws.Range("M132").Value = -589.5715999999998
